$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# --- Weekly crime table (rows 14-30) ---
# Row 14: Murder
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = $ws.Range("F14").NumberFormat
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 5
$ws.Range("L14").Value = -28.571428571428
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -86.842105263157

# Row 15: Rape
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = 21.428571428571
$ws.Range("L15").Value = 13.333333333333
$ws.Range("M15").Value = -15
$ws.Range("N15").Value = -32

# Row 16: Robbery
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = 15.384615384615
$ws.Range("F16").Value = 47
$ws.Range("G16").Value = 66
$ws.Range("H16").Value = -28.787878787878
$ws.Range("I16").Value = 180
$ws.Range("J16").Value = 232
$ws.Range("K16").Value = -22.413793103448
$ws.Range("L16").Value = 45.16129032258
$ws.Range("M16").Value = -4.761904761904
$ws.Range("N16").Value = -77.941176470588

# Row 17: Fel. Assault
$ws.Range("C17").Value = 27
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = 22.727272727272
$ws.Range("F17").Value = 95
$ws.Range("G17").Value = 90
$ws.Range("H17").Value = 5.555555555555
$ws.Range("I17").Value = 292
$ws.Range("J17").Value = 323
$ws.Range("K17").Value = -9.597523219814
$ws.Range("L17").Value = 39.712918660287
$ws.Range("M17").Value = 36.448598130841
$ws.Range("N17").Value = -25.699745547073

# Row 18: Burglary
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 40.90909090909
$ws.Range("I18").Value = 128
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = 4.065040650406
$ws.Range("L18").Value = 24.271844660194
$ws.Range("M18").Value = 6.666666666666
$ws.Range("N18").Value = -74.348697394789

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = 21.739130434782
$ws.Range("F19").Value = 87
$ws.Range("G19").Value = 95
$ws.Range("H19").Value = -8.421052631578
$ws.Range("I19").Value = 308
$ws.Range("J19").Value = 370
$ws.Range("K19").Value = -16.756756756756
$ws.Range("L19").Value = 8.833922261484
$ws.Range("M19").Value = 65.591397849462
$ws.Range("N19").Value = 36.283185840708

# Row 20: G.L.A.
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 83.333333333333
$ws.Range("F20").Value = 47
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 38.235294117647
$ws.Range("I20").Value = 169
$ws.Range("J20").Value = 157
$ws.Range("K20").Value = 7.64331210191
$ws.Range("L20").Value = 29.007633587786
$ws.Range("M20").Value = 141.428571428571
$ws.Range("N20").Value = -80.52995391705

# Row 21: TOTAL
$ws.Range("C21").Value = 88
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = 25.714285714285
$ws.Range("F21").Value = 310
$ws.Range("G21").Value = 310
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1099
$ws.Range("J21").Value = 1224
$ws.Range("K21").Value = -10.212418300653
$ws.Range("L21").Value = 26.032110091743
$ws.Range("M21").Value = 35.84672435105
$ws.Range("N21").Value = -61.640488656195

# Row 22: Transit
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 36
$ws.Range("K22").Value = -58.333333333333
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -21.052631578947

# Row 23: Housing
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 119
$ws.Range("J23").Value = 98
$ws.Range("K23").Value = 21.428571428571
$ws.Range("L23").Value = 58.666666666666
$ws.Range("M23").Value = 170.454545454545

# Row 24: Petit Larceny
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 63
$ws.Range("E24").Value = -53.968253968254
$ws.Range("F24").Value = 173
$ws.Range("G24").Value = 251
$ws.Range("H24").Value = -31.075697211155
$ws.Range("I24").Value = 656
$ws.Range("J24").Value = 830
$ws.Range("K24").Value = -20.963855421686
$ws.Range("L24").Value = 1.234567901234
$ws.Range("M24").Value = 48.41628959276

# Row 25: Misd. Assault
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = -58.823529411764
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = -6.930693069306
$ws.Range("I25").Value = 323
$ws.Range("J25").Value = 365
$ws.Range("K25").Value = -11.506849315068
$ws.Range("L25").Value = 31.836734693877
$ws.Range("M25").Value = -36.039603960396

# Row 26: UCR Rape*
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = $ws.Range("F26").NumberFormat
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 23
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -23.333333333333

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -15.151515151515
$ws.Range("L27").Value = -3.448275862068

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 25
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 19.047619047619
$ws.Range("L28").Value = 19.047619047619
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -78.260869565217

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 8
$ws.Range("H29").Value = -11.111111111111
$ws.Range("I29").Value = 21
$ws.Range("J29").Value = 19
$ws.Range("K29").Value = 10.526315789473
$ws.Range("L29").Value = 5
$ws.Range("M29").Value = 10.526315789473
$ws.Range("N29").Value = -80.188679245283
